# The post at old row 106 ("「誰がこれをしたの？」") was removed from the
# sheet. Deleting the entire row shifts every subsequent row (old 107..194)
# up by one (new 106..193), which matches the target diff, and Excel
# automatically shrinks the sheet's used-range dimension from A1:C194 to
# A1:C193.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(106).Delete()
